# Adds two blank paragraphs and a new "property description" paragraph
# at the end of the document, formatted the way Word/Outlook-pasted text
# from an email client ("normaltextrun" run style) looks, and registers
# the supporting custom character style in styles.xml.

$d = $word.ActiveDocument

# Register the custom character style used by the new run, based on
# Word's built-in "Default Paragraph Font" character style.
$style = $d.Styles.Add("normaltextrun", 2)
$style.BaseStyle = "DefaultParagraphFont"

# Build the new content: two empty paragraphs followed by a paragraph
# holding the new text, with the desired direct run formatting already
# baked into the run properties.
$newText = "- a lot of progress with little last if electrical plumbing door knobs towel rods etc. new counter depth fridge delivered. Should be on the market very soon."

$wordml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' `
        + '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' `
        + '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        +   '<w:r>' `
        +     '<w:rPr>' `
        +       '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' `
        +       '<w:color w:val="000000"/>' `
        +       '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' `
        +     '</w:rPr>' `
        +     '<w:t>' + $newText + '</w:t>' `
        +   '</w:r>' `
        + '</w:p>'

$insertionPoint = $d.Range($d.Content.End, $d.Content.End)
$insertionPoint.InsertXML($wordml)

# Apply the custom "normaltextrun" character style to the run that holds
# the new text (InsertXML does not understand <w:rStyle>, so it has to be
# applied afterwards through the object model).
$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)
$newRunRange = $d.Range($lastParagraph.Range.Start, $lastParagraph.Range.End - 1)
$newRunRange.Style = "normaltextrun"
